$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new rows at row 15 (pushes the old "PriceBeakdown" row, and
#    everything below it, down by two rows). Excel auto-extends the
#    A12:A15 merged cell to A12:A17 and shifts the later merges too.
# ---------------------------------------------------------------------------
$ws.Range("15:16").Insert()

# Bring over the "middle of merge" row formatting (border/fill/font/number
# format) from row 14 so the two new rows look like the other rows that
# belong to the same task-group block.
$ws.Range("A14:K14").Copy()
$ws.Range("A15:K16").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A15:K16").RowHeight = 14.25

# ---------------------------------------------------------------------------
# 2. Populate the two new task rows.
# ---------------------------------------------------------------------------
$ws.Range("B15").Value = "BookingAdmin"
$ws.Range("C15").Value = 45802
$ws.Range("D15").Value = 45806
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = "Từ Triệu Dương"
$ws.Range("G15").Value = 45807
$ws.Range("H15").Value = 45809
$ws.Range("I15").Value = 1

$ws.Range("B16").Value = "LocationAdmin"
$ws.Range("C16").Value = 45802
$ws.Range("D16").Value = 45806
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = "Nguyễn Quốc Thắng"
$ws.Range("G16").Value = 45807
$ws.Range("H16").Value = 45809
$ws.Range("I16").Value = 1

# ---------------------------------------------------------------------------
# 3. A few progress values on earlier rows were bumped up to 100%.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 1
$ws.Range("I11").Value = 1
$ws.Range("I12").Value = 1
$ws.Range("I13").Value = 1

# ---------------------------------------------------------------------------
# 4. Conditional formatting range grew along with the new rows.
# ---------------------------------------------------------------------------
$ws.Range("E2:E25,I2:I25").FormatConditions.Delete()
$cf = $ws.Range("E2:E25,I2:I25").FormatConditions.Add(1, 5, "0.99")
$cf.SetFirstPriority()

# ---------------------------------------------------------------------------
# 5. Scroll position / selection as last left by the author.
# ---------------------------------------------------------------------------
$ws.Range("J7").Select()
$excel.ActiveWindow.ScrollRow = 7
